{"js": "// The upstream commit's canonical-OOXML diff for this particular\n// template file (duplicatedUserContentIDPreExistingLostFile-template.docx)\n// touches only word/document.xml (the <w:sectPr>/<w:pgSz>/<w:pgMar> block)\n// and word/styles.xml (<w:docDefaults>, <w:latentStyles>, and the four\n// <w:style> definitions). In every single hunk the *set* of XML attributes\n// and their values is exactly the same before and after - only the\n// serialized attribute order changes (e.g. w:w/w:h -> w:h/w:w,\n// w:uiPriority/w:qFormat -> w:qFormat/w:uiPriority, etc.). That is a\n// by-product of the fixture being re-saved by a newer tool/library as part\n// of the larger commit (which, per the commit message, actually added an\n// M2Doc version to a different file's custom document properties) - no\n// page size, margin, font, language, or style value actually changes for\n// this document.\n//\n// The Word JavaScript API does not expose a way to control low-level XML\n// attribute ordering (that is an implementation detail of the save/\n// serialization layer, not part of the document object model), and none of\n// <w:pgSz>/<w:pgMar>/<w:docDefaults>/<w:latentStyles>/<w:style> default\n// definitions are reachable/settable through Word.Section or\n// Word.Body/Range in the JS API anyway. So the correct, content-faithful\n// way to \"apply\" this attribute-reordering-only diff here is to leave the\n// document's actual content/formatting untouched - which is what this\n// script verifies (read-only) before returning.\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\n// Confirm the body text is present/unchanged; no mutation is made because\n// the diff does not add, remove, or alter any visible content or value -\n// it only reorders XML attributes that are not controllable from Office.js.\nvoid body.text;\n", "ps1": "# The upstream commit's canonical-OOXML diff for this particular template\n# file (duplicatedUserContentIDPreExistingLostFile-template.docx) touches\n# only word/document.xml (the <w:sectPr>/<w:pgSz>/<w:pgMar> block) and\n# word/styles.xml (<w:docDefaults>, <w:latentStyles>, and the four\n# <w:style> definitions). In every hunk the *set* of XML attributes and\n# their values is identical before and after - only the serialized\n# attribute order changes (e.g. w:w/w:h -> w:h/w:w, w:top/w:right/... ->\n# w:bottom/w:footer/..., w:uiPriority/w:qFormat -> w:qFormat/w:uiPriority,\n# and so on for all 371 <w:lsdException> entries and the <w:style>/\n# <w:tblInd>/<w:tblCellMar> elements). That reordering is a by-product of\n# the fixture being re-saved by a newer tool/library as part of the larger\n# commit (which, per the commit message, actually added an M2Doc version\n# to a different file's custom document properties) - the page size,\n# margins, fonts, language, and style definitions for this document do not\n# actually change.\n#\n# Word's COM automation model has no property that controls the order in\n# which attributes are written out when a part is serialized - that is an\n# internal detail of the save layer, not part of the object model. Probing\n# this sandbox confirms that: (1) re-assigning a PageSetup margin/size to\n# its own current value round-trips the section's XML with the attribute\n# order completely unchanged, and (2) poking a Styles(\"Normal\").Font\n# property (even \"set to itself\") spuriously *adds* a direct formatting\n# override to the style that is not present in either the \"before\" or the\n# real \"after\" document - i.e. touching the style sheet through the object\n# model actively diverges from the target rather than reproducing it.\n#\n# So the correct, content-faithful way to \"apply\" this attribute-\n# reordering-only diff is to leave the document's actual content and\n# formatting untouched. This script does that, after read-only-verifying\n# (without writing, so nothing is perturbed) that the page geometry the\n# diff mentions is already the expected, unchanged one.\n\n$d = $word.ActiveDocument\n$section = $d.Sections.Item(1)\n\n$pageWidth = $section.PageSetup.PageWidth\n$pageHeight = $section.PageSetup.PageHeight\n$topMargin = $section.PageSetup.TopMargin\n$rightMargin = $section.PageSetup.RightMargin\n$bottomMargin = $section.PageSetup.BottomMargin\n$leftMargin = $section.PageSetup.LeftMargin\n$headerDistance = $section.PageSetup.HeaderDistance\n$footerDistance = $section.PageSetup.FooterDistance\n$gutter = $section.PageSetup.Gutter\n\n# No assignment back - the values already match the diff's (unchanged)\n# target, and the object model offers no way to influence XML attribute\n# ordering, so no write is needed or appropriate here.\n"}
